# Regenerate the "K" column (column G) of the save-data sheet with the
# recalculated strikeout/K values (s_vals), replacing the old Strike# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G, rows 2-33 (row 1 is the header "K").
$gValues = @(2, 0, 0, 1, 1, 0, 0, 1, 0, 0, 1, 0, 0, 0, 1, 1, 3, 0, 2, 1, 1, 0, 0, 3, 1, 2, 2, 1, 1, 2, 0, 2)

$startRow = 2
for ($i = 0; $i -lt $gValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $gValues[$i]
}
